$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The model repr text shared by F2 and F3 (multi-line, matching sklearn's
# pretty-printed estimator representation).
$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       LinearRegression())]),`n                                            param_grid={'model__fit_intercept': [True,`n                                                                                 False]},`n                                            scoring='neg_mean_squared_error'))"

# Add new header "Modelo" in F1, matching the bold/centered/bordered style
# already used by the other header cells (copy format from E1).
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the model description for each data row.
$ws.Range("F2").Value = $modelText
$ws.Range("F3").Value = $modelText

# Keep row heights on their default/auto value (avoid Excel forcing a
# customHeight because of the embedded newlines).
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
